$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.061319947242737
$ws.Range("B1").Value = 3.684934139251709
$ws.Range("C1").Value = 3.250450849533081
$ws.Range("D1").Value = 2.02033519744873
$ws.Range("E1").Value = 1.159409284591675
